$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New row of data for 22/03/2020 (row 21), matching the style/format of the
# previous data row (row 20): column A uses the date format, the rest use
# the default "General" style already applied on that row.
$ws.Range("A21").Value = 43912
$ws.Range("A21").NumberFormat = $ws.Range("A20").NumberFormat

$values = @{
    "B21" = 20
    "C21" = 2
    "D21" = 1
    "E21" = 11
    "F21" = 1
    "G21" = 6
    "H21" = 19
    "I21" = 459
    "J21" = 8
    "K21" = 28
    "L21" = 81
    "M21" = 46
    "N21" = 40
    "O21" = 4
    "P21" = 33
    "Q21" = 1
    "R21" = 6
    "S21" = 746
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# Update the active selection to N1, matching the diff's sheetView change.
$ws.Range("N1").Select()
